# "update database and change read_price algorithm"
#
# The source data for the quarter ending 1400/09 (spreadsheet column J) was
# recomputed. Previously most of these cells held the literal placeholder
# "-" (a shared string); the refreshed read_price algorithm now fills them
# in with real figures, and a handful of already-numeric cells are updated
# to their corrected totals.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Overview")

# --- مقدار فروش داخلی / مقدار فروش خارجی (quantities) ---
$ws.Range("J11").Value = 82839
$ws.Range("J12").Value = 184479
$ws.Range("J13").Value = 25558
$ws.Range("J14").Value = 201478
$ws.Range("J16").Value = 494354
$ws.Range("J18").Value = 230090
$ws.Range("J19").Value = 17649
$ws.Range("J20").Value = 198719
$ws.Range("J22").Value = 446458
$ws.Range("J26").Value = 0
$ws.Range("J27").Value = 940812

# --- مبلغ فروش داخلی / خارجی (sales amounts) ---
$ws.Range("J34").Value = 2923895
$ws.Range("J35").Value = 28821668
$ws.Range("J36").Value = 2161912
$ws.Range("J37").Value = 26938654
$ws.Range("J39").Value = 60846129
$ws.Range("J41").Value = 36879263
$ws.Range("J42").Value = 1409712
$ws.Range("J43").Value = 29571612
$ws.Range("J45").Value = 67860587
$ws.Range("J49").Value = 0
$ws.Range("J50").Value = 128706716

# --- نرخ فروش خارجی (unit price correction) ---
$ws.Range("J59").Value = 84588465

# --- بهای تمام شده داخلی / خارجی (cost of sales) ---
$ws.Range("J73").Value = -3099511
$ws.Range("J74").Value = -34923642
$ws.Range("J75").Value = -2109625
$ws.Range("J76").Value = -45226009
$ws.Range("J78").Value = -85358787
$ws.Range("J80").Value = -61949161
$ws.Range("J82").Value = -50830386
$ws.Range("J84").Value = -113648736
$ws.Range("J89").Value = -199007523

# --- سود ناخالص داخلی / خارجی (gross profit) ---
$ws.Range("J96").Value = 727357
$ws.Range("J97").Value = 30892467
$ws.Range("J98").Value = 1479588
$ws.Range("J99").Value = 22565977
$ws.Range("J101").Value = 55665389
$ws.Range("J103").Value = 47906155
$ws.Range("J105").Value = 27470939
$ws.Range("J107").Value = 75917617
$ws.Range("J110").Value = 131583006
